$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header column C1 from "audioFalse" to "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# Replace the two distinct audio file references in column C (rows 2-3)
# with the single shared value "train2P2"
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
